$d = $word.ActiveDocument

# Update the date/title line (unique text, safe to use Find & Replace).
$d.Content.Find.Execute("2023-07-20 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-21 Friday", 2) | Out-Null

# Update each table cell directly by (row, column) so that duplicate/overlapping
# values (e.g. "87" divided problems, or a new value equal to another cell's old
# value) cannot cause a global Find/Replace to touch the wrong cell.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "14÷9=1, 5"  # was "82÷7=11, 5"
$t.Cell(1,2).Range.Text = "37÷9=4, 1"  # was "32÷7=4, 4"
$t.Cell(1,3).Range.Text = "80÷5=16, 0"  # was "31÷8=3, 7"
$t.Cell(1,4).Range.Text = "91÷9=10, 1"  # was "39÷7=5, 4"
$t.Cell(1,5).Range.Text = "53÷9=5, 8"  # was "77÷7=11, 0"
$t.Cell(5,1).Range.Text = "64÷4=16, 0"  # was "49÷9=5, 4"
$t.Cell(5,2).Range.Text = "44÷3=14, 2"  # was "76÷9=8, 4"
$t.Cell(5,3).Range.Text = "63÷9=7, 0"  # was "38÷2=19, 0"
$t.Cell(5,4).Range.Text = "40÷4=10, 0"  # was "48÷5=9, 3"
$t.Cell(5,5).Range.Text = "79÷8=9, 7"  # was "98÷7=14, 0"
$t.Cell(9,1).Range.Text = "91÷7=13, 0"  # was "79÷8=9, 7"
$t.Cell(9,2).Range.Text = "70÷8=8, 6"  # was "23÷4=5, 3"
$t.Cell(9,3).Range.Text = "84÷2=42, 0"  # was "96÷9=10, 6"
$t.Cell(9,4).Range.Text = "90÷9=10, 0"  # was "56÷8=7, 0"
$t.Cell(9,5).Range.Text = "27÷9=3, 0"  # was "92÷5=18, 2"
$t.Cell(13,1).Range.Text = "39÷6=6, 3"  # was "52÷7=7, 3"
$t.Cell(13,2).Range.Text = "65÷2=32, 1"  # was "79÷4=19, 3"
$t.Cell(13,3).Range.Text = "39÷7=5, 4"  # was "90÷7=12, 6"
$t.Cell(13,4).Range.Text = "33÷9=3, 6"  # was "86÷8=10, 6"
$t.Cell(13,5).Range.Text = "94÷7=13, 3"  # was "87÷4=21, 3"
$t.Cell(17,1).Range.Text = "12÷8=1, 4"  # was "13÷7=1, 6"
$t.Cell(17,2).Range.Text = "47÷5=9, 2"  # was "25÷9=2, 7"
$t.Cell(17,3).Range.Text = "42÷8=5, 2"  # was "53÷2=26, 1"
$t.Cell(17,4).Range.Text = "42÷5=8, 2"  # was "62÷8=7, 6"
$t.Cell(17,5).Range.Text = "35÷8=4, 3"  # was "87÷4=21, 3"
